$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Block: rows 10-17 ---
$ws.Range("J10").Value = 1
$ws.Range("M10").Value = ""
$ws.Range("J11").Value = 1
$ws.Range("M12").Value = ""
$ws.Range("M14").Value = "Fly Ball"
$ws.Range("M15").Value = "Out"
$ws.Range("J17").Value = "FB,CB,CH"

# --- Block: rows 19-26 ---
$ws.Range("J19").Value = 3
$ws.Range("M19").Value = ""
$ws.Range("J20").Value = 2
$ws.Range("M21").Value = ""
$ws.Range("J23").Value = "Roblez"
$ws.Range("M23").Value = ""
$ws.Range("M24").Value = "Undefined"
$ws.Range("J25").Value = "88-90 MPH"
$ws.Range("J26").Value = "FB,CB,CH"

# --- Block: rows 28-35 ---
$ws.Range("M28").Value = ""
$ws.Range("M30").Value = ""
$ws.Range("M32").Value = ""
$ws.Range("J35").Value = "SL,FB,CB,CH"

# --- Block: rows 37-44 ---
$ws.Range("J37").Value = 6
$ws.Range("M37").Value = ""
$ws.Range("J38").Value = 0
$ws.Range("M39").Value = ""
$ws.Range("J41").Value = "Herbst"
$ws.Range("M41").Value = ""
$ws.Range("J42").Value = "Right"
$ws.Range("M42").Value = "Undefined"
$ws.Range("J43").Value = "83-85 MPH"
$ws.Range("J44").Value = "SL,FB,CB,CH"

# --- Block: rows 46-53 ---
$ws.Range("J46").Value = 7
$ws.Range("M46").Value = ""
$ws.Range("M48").Value = ""
$ws.Range("J50").Value = "Plum"
$ws.Range("M50").Value = "Line Drive"
$ws.Range("M51").Value = "Single"
$ws.Range("J52").Value = "84-86 MPH"
$ws.Range("J53").Value = "SL,FB,CH"

# --- Block: rows 61-68 ---
$ws.Range("J61").Value = 8
$ws.Range("M61").Value = ""
$ws.Range("J62").Value = 1
$ws.Range("M63").Value = ""
$ws.Range("J65").Value = "Thompson"
$ws.Range("M65").Value = "Ground Ball"
$ws.Range("J66").Value = "Left"
$ws.Range("M66").Value = "Out"
$ws.Range("J67").Value = "84-84 MPH"
$ws.Range("J68").Value = "SL,FB,CH"
